$d = $word.ActiveDocument

# The "АТИ: ID:" label run loses its text (the run itself, and the line
# break before it, stay in place) and the placeholder that followed it
# is renamed from {{EXECUTOR_ID}} to {{EXECUTOR_CPP}}.

$found1 = $d.Content.Find.Execute("АТИ: ID: ", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 2)

$found2 = $d.Content.Find.Execute("{{EXECUTOR_ID}}", $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "{{EXECUTOR_CPP}}", 2)

Write-Output "removed label: $found1; renamed placeholder: $found2"
